$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.033.48"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "3.385.88"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.42%  "
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.387"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("D12").Value = "3.964.63"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000171"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.378.05"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "61.130.56"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "385.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.19%  "
$ws.Range("E23").Value = "  -1.65%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  -1.58%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.185"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.39%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.51%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.99%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.08%  "
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "166.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.419.67"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.98%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0768"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.68%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.777"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.455.54"
$ws.Range("E46").Value = "  -2.94%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.82%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0263"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.06%  "
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.205"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.42%  "
